# Update the header row (row 1) to show the meeting day (Tuesday by default)
# instead of the week date range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "'02/01/2024"
$ws.Range("B1").Value = "'09/01/2024"
$ws.Range("C1").Value = "'16/01/2024"
$ws.Range("D1").Value = "'23/01/2024"
$ws.Range("E1").Value = "'30/01/2024"
$ws.Range("F1").Value = "'06/02/2024"
$ws.Range("G1").Value = "'13/02/2024"
$ws.Range("H1").Value = "'20/02/2024"
$ws.Range("I1").Value = "'27/02/2024"

# Append the opening song/prayer + introductory words text to row 3.
$ws.Range("A3").Value = $ws.Range("A3").Value + "y oración | Palabras de introducción(1 min.)"
$ws.Range("B3").Value = $ws.Range("B3").Value + "y oración | Palabras de introducción(1 min.)"
$ws.Range("C3").Value = $ws.Range("C3").Value + "y oración | Palabras de introducción(1 min.)"
$ws.Range("D3").Value = $ws.Range("D3").Value + "y oración | Palabras de introducción(1 min.)"
$ws.Range("E3").Value = $ws.Range("E3").Value + "y oración | Palabras de introducción(1 min.)"
$ws.Range("F3").Value = $ws.Range("F3").Value + "y oración | Palabras de introducción(1 min.)"
$ws.Range("G3").Value = $ws.Range("G3").Value + "y oración | Palabras de introducción(1 min.)"
$ws.Range("H3").Value = $ws.Range("H3").Value + "y oración | Palabras de introducción(1 min.)"
$ws.Range("I3").Value = $ws.Range("I3").Value + "y oración | Palabras de introducción(1 min.)"

# Fill in the "NUESTRA VIDA CRISTIANA" section header that now starts one row
# earlier (row 11 instead of row 13), shifting the whole program up by two rows.
$ws.Range("A11").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("F11").Value = "NUESTRA VIDA CRISTIANA"

$ws.Range("A12").Value = "Canción 116"
$ws.Range("B12").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("C12").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("D12").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("E12").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("F12").Value = "Canción 32"
$ws.Range("H12").Value = "NUESTRA VIDA CRISTIANA"

$ws.Range("A13").Value = "6. Necesidades de la congregación"
$ws.Range("B13").Value = "Canción 58"
$ws.Range("C13").Value = "Canción 49"
$ws.Range("D13").Value = "Canción 111"
$ws.Range("E13").Value = "Canción 108"
$ws.Range("F13").Value = "6. Necesidades de la congregación"
$ws.Range("H13").Value = "Canción 10"

$ws.Range("A14").Value = "7. Estudio bíblico de la congregación"
$ws.Range("B14").Value = "7. ¿“Predica la palabra” informalmente con entusiasmo?"
$ws.Range("C14").Value = "7. Esté preparado por si necesita tratamiento médico o una intervención quirúrgica"
$ws.Range("D14").Value = "7. Meditar en la creación nos ayuda a no perder de vista el cuadro completo"
$ws.Range("E14").Value = "7. Ayude a los demás a sentir el amor de Jehová"
$ws.Range("F14").Value = "7. Estudio bíblico de la congregación"
$ws.Range("H14").Value = "7. Cómo predicar informalmente de forma natural"

$ws.Range("A15").Value = "Palabras de conclusión(3 mins.)|Canción 54y oración"
$ws.Range("B15").Value = "8. Estudio bíblico de la congregación"
$ws.Range("C15").Value = "8. Estudio bíblico de la congregación"
$ws.Range("D15").Value = "8. Estudio bíblico de la congregación"
$ws.Range("E15").Value = "¿Cómo podemos ayudar a nuestros hermanos a creer que Jehová los ama?"
$ws.Range("F15").Value = "Palabras de conclusión(3 mins.)|Canción 61y oración"
$ws.Range("H15").Value = "Las siguientes ideas le servirán para comenzar conversaciones:"

# The old rows 16-18 are no longer needed; their content was folded into
# rows 11-15 above, so drop the now-empty trailing rows.
$ws.Rows("16:18").Delete()
